$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly rotation of the Fruta/Hortaliza price data.
# Each data row (2-8) receives the Fecha/Volumen/Precio values that
# previously belonged to another row, per the new weekly reporting cycle.

$ws.Range("D2").Value = 44382
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7438
$ws.Range("P2").Value = 124

$ws.Range("D3").Value = 44421
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("P3").Value = 142

$ws.Range("D4").Value = 44362
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("P4").Value = 142

$ws.Range("D5").Value = 44400
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9500
$ws.Range("P5").Value = 158

$ws.Range("D6").Value = 44281
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 5500
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5750
$ws.Range("P6").Value = 96

$ws.Range("D7").Value = 44242
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = 5250
$ws.Range("P7").Value = 88

$ws.Range("D8").Value = 44494
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("P8").Value = 92
